# "added aura to legendaries"
# Rename the nicknames (Column1) of four of the Pokemon in the supply table
# (Sheet1 / Table1). The formulas in columns I and J recompute automatically
# since they reference Table1[[#This Row],[Column1]].
#
# The shared-string table order in the saved file reflects the order in
# which new distinct strings are introduced, so the cells are written in
# the same row order (C7, C9, C5, C6) that reproduces the canonical
# ordering of the newly-added strings (quartzan, azrure, diwa, milid).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7").Value = "quartzan"   # Articuno: articane -> quartzan
$ws.Range("C9").Value = "azrure"     # Moltres:  moltyte  -> azrure
$ws.Range("C5").Value = "diwa"       # Jirachi:  jiratik  -> diwa
$ws.Range("C6").Value = "milid"      # Darkrai:  draco    -> milid

$ws.Activate()
$ws.Range("C7").Select() | Out-Null
